# Update vm_pu results for the 380 kV case: bus 0's voltage setpoint
# moves from 1.05 pu to 1.02 pu, which the underlying power-flow
# solver then propagates into every other bus voltage in the table
# (rows 2-25, columns B-F and I-N; columns A, G and H are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ B=1.02; C=1.017487223972259; D=1.02931213346166; E=1.018881255138838; F=1.015850058452286; I=1.03127810702048; J=1.022701060423484; K=1.032126410341878; L=1.021726154768394; M=1.018703985416108; N=1.011669028973444 }
    3 = @{ B=1.02; C=1.01835376443172; D=1.029776817588451; E=1.019613570787166; F=1.017357151644581; I=1.031429782335551; J=1.02320394480444; K=1.032400138692849; L=1.022264578470865; M=1.020014387534297; N=1.011836079773503 }
    4 = @{ B=1.02; C=1.018914469661225; D=1.030077299771513; E=1.020087827955235; F=1.018332298454257; I=1.031526382426304; J=1.023528726554068; K=1.032576335224361; L=1.022612710326617; M=1.020861779680625; N=1.01194393703005 }
    5 = @{ B=1.02; C=1.019150188680727; D=1.030203573104563; E=1.02028730059287; F=1.018742244339011; I=1.031566622907646; J=1.0236651165352; K=1.032650186031522; L=1.022759001212207; M=1.021217901296118; N=1.011989223633202 }
    6 = @{ B=1.02; C=1.019189766843842; D=1.030224772010273; E=1.020320798468572; F=1.01881107584612; I=1.031573357736513; J=1.023688008309982; K=1.03266257284691; L=1.022783560343519; M=1.021277688635989; N=1.011996824130689 }
    7 = @{ B=1.02; C=1.01891761935655; D=1.030078987238034; E=1.020090492946452; F=1.018337776188531; I=1.031526921576921; J=1.023530549586644; K=1.032577322896888; L=1.022614665323681; M=1.020866538669908; N=1.01194454237439 }
    8 = @{ B=1.02; C=1.017780076043553; D=1.029469215984238; E=1.01912866132626; F=1.016359399335016; I=1.031329685756408; J=1.022871140352561; K=1.032219108841939; L=1.021908172016391; M=1.019146953074435; N=1.011725533309507 }
    9 = @{ B=1.02; C=1.015775563108537; D=1.028393263782479; E=1.017436890588403; F=1.0128727125858; I=1.030970330652424; J=1.021704455974072; K=1.031580857031715; L=1.020661236985051; M=1.016112627444113; N=1.011337811583503 }
    10 = @{ B=1.02; C=1.01443923093143; D=1.027675077339633; E=1.016311174283142; F=1.010547610988209; I=1.030722859125941; J=1.020923507278758; K=1.031150685971699; L=1.019828625830677; M=1.014086682790126; N=1.011078127862708 }
    11 = @{ B=1.02; C=1.013860590171416; D=1.02736390463972; E=1.015824241447203; F=1.009540599817928; I=1.030613833683442; J=1.020584601112025; K=1.030963321697168; L=1.019467787566779; M=1.013208646544352; N=1.010965398085158 }
    12 = @{ B=1.02; C=1.013645657395943; D=1.02724829380366; E=1.015643450192308; F=1.009166511384362; I=1.030573056523672; J=1.020458603685779; K=1.030893562349268; L=1.019333709605487; M=1.012882381557092; N=1.010923482492095 }
    13 = @{ B=1.02; C=1.013691761195465; D=1.027273093936771; E=1.015682227021974; F=1.00924675650702; I=1.03058181603729; J=1.020485635654205; K=1.030908533368429; L=1.019362471885179; M=1.012952372156785; N=1.010932475461566 }
    14 = @{ B=1.02; C=1.013842823756706; D=1.027354348762568; E=1.015809295605169; F=1.009509678392401; I=1.03061046874504; J=1.020574188421478; K=1.030957558706021; L=1.019456705588666; M=1.01318167991995; N=1.010961934200442 }
    15 = @{ B=1.02; C=1.01393589841386; D=1.027404408953737; E=1.01588759701424; F=1.009671667822847; I=1.03062808550175; J=1.020628733788942; K=1.030987743142735; L=1.019514759923142; M=1.013322947530232; N=1.010980079041521 }
    16 = @{ B=1.02; C=1.014477633426438; D=1.027695724911785; E=1.016343501211346; F=1.010614437749121; I=1.03073005545699; J=1.020945983569446; K=1.0311630976828; L=1.019852566929909; M=1.014144938141879; N=1.011085603369246 }
    17 = @{ B=1.02; C=1.014817449589588; D=1.027878409109007; E=1.016629614839232; F=1.011205747763272; I=1.030793518572678; J=1.021144785209788; K=1.031272799832006; L=1.020064381008362; M=1.014660336618111; N=1.011151719747342 }
    18 = @{ B=1.02; C=1.015015658585884; D=1.027984946960511; E=1.016796549215266; F=1.011550627378682; I=1.030830355205504; J=1.021260670448846; K=1.03133668129616; L=1.02018789843791; M=1.014960884035016; N=1.011190256806886 }
    19 = @{ B=1.02; C=1.015083242721969; D=1.028021270385804; E=1.016853477836779; F=1.011668218935366; I=1.03084288494361; J=1.021300172037322; K=1.031358445225631; L=1.020230009567149; M=1.015063350224497; N=1.011203392280303 }
    20 = @{ B=1.02; C=1.014780990532145; D=1.027858810751427; E=1.016598912467068; F=1.01114230805951; I=1.030786728229685; J=1.021123463162283; K=1.031261040780918; L=1.020041658477248; M=1.014605047112591; N=1.011144628930033 }
    21 = @{ B=1.02; C=1.013798339591511; D=1.02733042198637; E=1.015771874935634; F=1.009432255589612; I=1.030602038969318; J=1.020548114945576; K=1.030943126483237; L=1.019428957394461; M=1.013114157938795; N=1.010953260515087 }
    22 = @{ B=1.02; C=1.013180508995158; D=1.02699804501512; E=1.015252331077349; F=1.008356845567309; I=1.030484296263105; J=1.020185719355757; K=1.030742293256814; L=1.019043458661276; M=1.012176061416318; N=1.010832692389759 }
    23 = @{ B=1.02; C=1.013508032524834; D=1.027174258783478; E=1.015527708462222; F=1.008926964639655; I=1.030546867413042; J=1.020377893812617; K=1.030848848282859; L=1.019247844245998; M=1.012673433517302; N=1.01089663123874 }
    24 = @{ B=1.02; C=1.014797464810807; D=1.02786766646373; E=1.01661278539566; F=1.011170973796057; I=1.030789797050603; J=1.021133097901357; K=1.03126635451738; L=1.020051925904059; M=1.014630030293273; N=1.01114783305012 }
    25 = @{ B=1.02; C=1.016293777521099; D=1.028671586101681; E=1.017873881703742; F=1.013774197225789; I=1.031064627606064; J=1.022006629743806; K=1.031746687782904; L=1.020983834312269; M=1.01689759523932; N=1.011438259443949 }
}

foreach ($row in $newValues.Keys) {
    $rowVals = $newValues[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
